$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.580.86'
$ws.Range('E2').Value = '  -1.09%  '

$ws.Range('D3').Value = '''1.858.02'
$ws.Range('E3').Value = '  -0.22%  '

$ws.Range('D4').Value = '''0.9997'
$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '''242.30'
$ws.Range('E5').Value = '  -0.96%  '

$ws.Range('D6').Value = '''0.6327'
$ws.Range('E6').Value = '  -3.97%  '

$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').Value = '''0.07601'
$ws.Range('E8').Value = '  +0.18%  '

$ws.Range('D9').Value = '''0.2994'
$ws.Range('E9').Value = '  -0.12%  '

$ws.Range('D10').Value = '''24.64'
$ws.Range('E10').Value = '  -0.50%  '

$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '''1.950.11'
$ws.Range('E11').Value = '  +4.55%  '

$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '''0.07719'
$ws.Range('E12').Value = '  +0.93%  '

$ws.Range('D13').Value = '''0.6933'
$ws.Range('E13').Value = '  -0.05%  '

$ws.Range('D14').Value = '''5.025'
$ws.Range('E14').Value = '  -1.18%  '

$ws.Range('D15').Value = '''83.68'
$ws.Range('E15').Value = '  -0.41%  '

$ws.Range('D16').Value = '''0.000009909'
$ws.Range('E16').Value = '  +1.94%  '

$ws.Range('D17').Value = '''2.150.34'
$ws.Range('E17').Value = '  +1.11%  '

$ws.Range('D18').Value = '''6.234'
$ws.Range('E18').Value = '  +1.37%  '

$ws.Range('D19').Value = '''29.674.36'
$ws.Range('E19').Value = '  -0.78%  '

$ws.Range('D20').Value = '''234.36'
$ws.Range('E20').Value = '  -1.14%  '

$ws.Range('D21').Value = '''12.60'
$ws.Range('E21').Value = '  -1.00%  '

$ws.Range('E22').Value = '  +0.02%  '

$ws.Range('D23').Value = '''7.678'
$ws.Range('E23').Value = '  -1.90%  '

$ws.Range('D24').Value = '''1.001'
$ws.Range('E24').Value = '  +0.11%  '

$ws.Range('D25').Value = '''155.63'
$ws.Range('E25').Value = '  -1.85%  '

$ws.Range('D26').Value = '''0.1401'
$ws.Range('E26').Value = '  -3.04%  '

$ws.Range('D27').Value = '''8.502'
$ws.Range('E27').Value = '  -1.37%  '

$ws.Range('D28').Value = '''17.76'
$ws.Range('E28').Value = '  -0.94%  '

$ws.Range('E29').Value = '  -1.01%  '

$ws.Range('D30').Value = '''0.05792'
$ws.Range('E30').Value = '  -4.62%  '

$ws.Range('E31').Value = '  -2.35%  '

$ws.Range('E32').Value = '  -0.72%  '

$ws.Range('D33').Value = '''4.042'
$ws.Range('E33').Value = '  -1.37%  '

$ws.Range('D34').Value = '''1.900'
$ws.Range('E34').Value = '  +1.18%  '

$ws.Range('D35').Value = '''1.173'
$ws.Range('E35').Value = '  -0.96%  '

$ws.Range('D36').Value = '''0.7257'
$ws.Range('E36').Value = '  -1.39%  '

$ws.Range('D37').Value = '''2.590'
$ws.Range('E37').Value = '  -0.75%  '

$ws.Range('D38').Value = '''1.256.91'
$ws.Range('E38').Value = '  +3.88%  '

$ws.Range('D39').Value = '''2.817'
$ws.Range('E39').Value = '  +0.04%  '

$ws.Range('D40').Value = '''0.01812'
$ws.Range('E40').Value = '  +0.70%  '

$ws.Range('D41').Value = '''0.9053'
$ws.Range('E41').Value = '  -1.12%  '

$ws.Range('D42').Value = '''6.151'
$ws.Range('E42').Value = '  -2.78%  '

$ws.Range('D43').Value = '''2.069.14'
$ws.Range('E43').Value = '  +1.61%  '

$ws.Range('E44').Value = '  -0.02%  '

$ws.Range('D45').Value = '''68.02'
$ws.Range('E45').Value = '  +0.66%  '

$ws.Range('D46').Value = '''101.59'
$ws.Range('E46').Value = '  -0.09%  '

$ws.Range('D47').Value = '''7.352'
$ws.Range('E47').Value = '  -2.73%  '

$ws.Range('D48').Value = '''0.00000000118'
$ws.Range('E48').Value = '  -2.53%  '

$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').Value = '''0.4053'
$ws.Range('E49').Value = '  -0.71%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''9.180'
$ws.Range('E50').Value = '  -0.18%  '

$ws.Range('D51').Value = '''1.717'
$ws.Range('E51').Value = '  +1.61%  '
